# Auto-generated PowerShell COM-interop script
# Applies Part 1 dialogue translation update (column C, English text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 58.25
$ws.Columns.Item(2).ColumnWidth = 64.42
$ws.Columns.Item(3).ColumnWidth = 61.59

# --- New column C (English translation) dialogue text + wrap styling ---
$ws.Range("C9").Value = "\n<Shina>Shit...`nI got beat, nyan."
$ws.Range("C9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30
$ws.Range("C10").Value = "\n<Lily>Shina!? What happened to you!"
$ws.Range("C10").WrapText = $false
$ws.Range("C11").Value = "\n<Lily>What do you think..."
$ws.Range("C11").WrapText = $false
$ws.Range("C12").Value = "\n<Lily>That form...`nKeheheheh♥"
$ws.Range("C12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 30
$ws.Range("C13").Value = "\n<Shina>Don't you laugh at me!!"
$ws.Range("C13").WrapText = $false
$ws.Range("C14").Value = "\n<Lily>Aw, kitty kitty♥`nPoor little thing."
$ws.Range("C14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 30
$ws.Range("C15").Value = "\n<Shina>I'll fucking kill you!!"
$ws.Range("C15").WrapText = $false
$ws.Range("C16").Value = "\n<Lily>Ow ow ow ow!`nI'm sorry, I'll stop, I'll stop!"
$ws.Range("C16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 30
$ws.Range("C17").Value = "\n<Shina>..."
$ws.Range("C17").WrapText = $false
$ws.Range("C18").Value = "\n<Shina>\n[1] made me cum, nya...`nShit."
$ws.Range("C18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 30
$ws.Range("C19").Value = "\n<Lily>\n[1]...`nI didn't give you enough credit."
$ws.Range("C19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 30
$ws.Range("C20").Value = "\n<Lily>First of all, I'm going to take care of that one.`nThen I'll put you back to normal.`nJust hold out for a little while, Ok?"
$ws.Range("C20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 45
$ws.Range("C21").Value = "\n<Shina>Put me back now, nya."
$ws.Range("C21").WrapText = $false
$ws.Range("C22").Value = "\n<Lily>I need to conserve my magic right now.`nJust wait a little while."
$ws.Range("C22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 30
$ws.Range("C23").Value = "\n<Shina>Tsk...`nMake it quick then, nya."
$ws.Range("C23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 30
$ws.Range("C24").Value = "\n<Shina>I let my guard down, nya.`nIf this was for real, I would've won.`nI'm so pissed, nya."
$ws.Range("C24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 45
$ws.Range("C25").Value = "\n<Lily>Of course not...`nYou would never lose to someone like that."
$ws.Range("C25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 30
$ws.Range("C26").Value = "\n<Shina>Lily, don't let your guard down too, nya."
$ws.Range("C26").WrapText = $false
$ws.Range("C27").Value = "\n<Lily>Right.`nWe're not playing around anymore."
$ws.Range("C27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 30
$ws.Range("C28").Value = "\n<Lily>I'll squeeze them dry, and turn it into magical power.`nI'll destroy \n[1] and avenge you.`nLeave it all to me."
$ws.Range("C28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 45
$ws.Range("C29").Value = "\n<Lily>I won't let my guard down... Unlike you...`nHeh heh heh♥"
$ws.Range("C29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 30
$ws.Range("C31").Value = "\n<Lime>Awaaaahhh..."
$ws.Range("C31").WrapText = $false
$ws.Range("C32").Value = "\n<Lily>L-Lime!?`nWhat happened!?`nYou're melting!!"
$ws.Range("C32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 45
$ws.Range("C33").Value = "\n<Lime>I was beaten..."
$ws.Range("C33").WrapText = $false
$ws.Range("C34").Value = "\n<Lily>\n[1] had sex with you...!?`nWhat in the world..."
$ws.Range("C34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 30
$ws.Range("C35").Value = "\n<Lime>But...`nIt...`nIt felt so good...♥"
$ws.Range("C35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 45
$ws.Range("C36").Value = "\n<Lily>Are you...`nAre you alright?`nDo you need some water?"
$ws.Range("C36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 45
$ws.Range("C37").Value = "\n<Lime>I'm OK-♥`nIt doesn't feel like I'll melt anymore.`nI can't take my succubus form right now though."
$ws.Range("C37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 45
$ws.Range("C38").Value = "\n<Lily>Ah...`nYour magical energy has been completely sapped away..."
$ws.Range("C38").WrapText = $true
$ws.Rows.Item(38).RowHeight = 30
$ws.Range("C39").Value = "\n<Lime>I lose my shape if I lose my magical power?"
$ws.Range("C39").WrapText = $false
$ws.Range("C40").Value = "\n<Lily>That's right.`nYour magic coats your body, holding you in that shape.`nWithout that power, you'll fall apart."
$ws.Range("C40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 60
$ws.Range("C41").Value = "\n<Lily>Simply put, you create a magical, humanoid vessel, and`nyour body is sort of poured into that image.`nThat's not completely accurate, but..."
$ws.Range("C41").WrapText = $true
$ws.Rows.Item(41).RowHeight = 45
$ws.Range("C42").Value = "\n<Lime>I get it."
$ws.Range("C42").WrapText = $false
$ws.Range("C43").Value = "\n<Lily>Everybody's vessel takes a different shape.`nYour natural shape is that of a slime.`nMy power transforms your vessel into a humanoid form."
$ws.Range("C43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 45
$ws.Range("C44").Value = "\n<Lily>I'll get you back to normal in a little while.`nCan you hold on for now?"
$ws.Range("C44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 30
$ws.Range("C45").Value = "\n<Lime>Sure.`nI don't mind taking this shape every once in a while.`nThe low viewpoint is quite interesting-♥"
$ws.Range("C45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 45
$ws.Range("C46").Value = "\n<Lily>\n[1]...`nNo matter how much you beg, or cry,`nor scream and cum... I'll never forgive you."
$ws.Range("C46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 45

# --- Column B wrap styling updates on rows 26-46 (matches new font/alignment applied in this edit) ---
$ws.Range("B26").WrapText = $false
$ws.Range("B29").WrapText = $true
$ws.Range("B31").WrapText = $false
$ws.Range("B32").WrapText = $true
$ws.Range("B33").WrapText = $false
$ws.Range("B34").WrapText = $true
$ws.Range("B35").WrapText = $true
$ws.Range("B36").WrapText = $true
$ws.Range("B37").WrapText = $true
$ws.Range("B38").WrapText = $true
$ws.Range("B39").WrapText = $false
$ws.Range("B40").WrapText = $true
$ws.Range("B41").WrapText = $true
$ws.Range("B42").WrapText = $false
$ws.Range("B43").WrapText = $true
$ws.Range("B44").WrapText = $true
$ws.Range("B45").WrapText = $true
$ws.Range("B46").WrapText = $true

# --- Selection / view state to match target (active cell C46) ---
$ws.Range("C46").Select()
